$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new value would otherwise be
# auto-detected as a number by Excel (mirrors typing into a text-formatted column).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '42.212.50'
$ws.Range("E2").Value = '  -0.59%  '
$ws.Range("D3").Value = '2.262.93'
$ws.Range("E3").Value = '  -0.77%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '307.90'
$ws.Range("E5").Value = '  +0.23%  '
$ws.Range("D6").Value = '97.23'
$ws.Range("E6").Value = '  -1.28%  '
$ws.Range("E7").Value = '  -0.83%  '
$ws.Range("D9").Value = '0.488'
$ws.Range("E9").Value = '  -1.61%  '
$ws.Range("D10").Value = '34.88'
$ws.Range("E10").Value = '  -3.46%  '
$ws.Range("D11").Value = '0.0806'
$ws.Range("E11").Value = '  +0.64%  '
$ws.Range("E12").Value = '  +1.06%  '
$ws.Range("D13").Value = '6.83'
$ws.Range("E13").Value = '  +1.59%  '
$ws.Range("D14").Value = '2.612.58'
$ws.Range("E14").Value = '  -0.82%  '
$ws.Range("D15").Value = '14.54'
$ws.Range("E15").Value = '  +0.36%  '
$ws.Range("D16").Value = '2.262.03'
$ws.Range("E16").Value = '  -1.15%  '
$ws.Range("D17").Value = '0.785'
$ws.Range("E17").Value = '  -1.79%  '
$ws.Range("D18").Value = '42.089.71'
$ws.Range("E18").Value = '  -0.58%  '
$ws.Range("D19").Value = '12.27'
$ws.Range("E19").Value = '  -3.22%  '
$ws.Range("D20").Value = '0.0₃0902'
$ws.Range("E20").Value = '  -1.49%  '
$ws.Range("D21").Value = '5.94'
$ws.Range("E21").Value = '  -1.38%  '
$ws.Range("D22").Value = '67.60'
$ws.Range("E22").Value = '  -0.19%  '
$ws.Range("D23").Value = '236.20'
$ws.Range("E23").Value = '  -2.48%  '
$ws.Range("D24").Value = '2.58'
$ws.Range("E24").Value = '  -0.90%  '
$ws.Range("D25").Value = '1.97'
$ws.Range("E25").Value = '  +0.62%  '
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  -0.14%  '
$ws.Range("D27").Value = '23.49'
$ws.Range("E27").Value = '  -1.76%  '
$ws.Range("D28").Value = '37.05'
$ws.Range("E28").Value = '  -2.33%  '
$ws.Range("D29").Value = '9.54'
$ws.Range("E29").Value = '  -0.24%  '
$ws.Range("E30").Value = '  +0.25%  '
$ws.Range("D31").Value = '163.08'
$ws.Range("E31").Value = '  +1.69%  '
$ws.Range("D32").Value = '5.22'
$ws.Range("E32").Value = '  -0.76%  '
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("D34").Value = '3.08'
$ws.Range("E34").Value = '  -2.00%  '
$ws.Range("D35").Value = '17.58'
$ws.Range("E35").Value = '  +2.09%  '
$ws.Range("D36").Value = '0.0731'
$ws.Range("E36").Value = '  -2.31%  '
$ws.Range("E37").Value = '  -0.03%  '
$ws.Range("E38").Value = '  -2.99%  '
$ws.Range("E39").Value = '  -0.33%  '
$ws.Range("D40").Value = '1.80'
$ws.Range("E40").Value = '  -2.63%  '
$ws.Range("D41").Value = '4.13'
$ws.Range("E41").Value = '  -0.92%  '
$ws.Range("E42").Value = '  -4.86%  '
$ws.Range("D43").Value = '1.946.17'
$ws.Range("E43").Value = '  -2.92%  '
$ws.Range("E44").Value = '  -1.85%  '
$ws.Range("D45").Value = '18.53'
$ws.Range("E45").Value = '  -3.52%  '
$ws.Range("D46").Value = '2.93'
$ws.Range("E46").Value = '  -2.76%  '
$ws.Range("E47").Value = '  -2.87%  '
$ws.Range("D48").Value = '54.42'
$ws.Range("E48").Value = '  +1.50%  '
$ws.Range("D49").Value = '2.488.53'
$ws.Range("E49").Value = '  -0.57%  '
$ws.Range("D50").Value = '91.62'
$ws.Range("E50").Value = '  -1.25%  '
$ws.Range("D51").Value = '71.37'
$ws.Range("E51").Value = '  -2.75%  '
